$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A31").Value = 46000
$ws.Range("B31").Value = 5

$ws.Range("A31:B31").Select() | Out-Null
